$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I (rows 2-27) currently holds a formula =FALSE() (boolean result).
# Replace each with the literal text "False" (a shared string, not a boolean
# and not a formula). A leading apostrophe forces Excel to store it as text.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 9).Value = "'False"
}

# New blank rows appended below the table (I28:I32), keeping the same style
# as the column above but with no value.
for ($r = 28; $r -le 32; $r++) {
    $ws.Cells.Item($r, 9).Value = ""
}

$ws.Range("I2:I32").NumberFormat = "@"

# Sheet view state: scrolled down, selection now on H28:J34.
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("H28:J34").Select()
